$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 65000
$ws.Range("J3").Value = 65000
$ws.Range("L3").Value = 65000
$ws.Range("N3").Value = -65228
$ws.Range("H32").Value = 638.2
$ws.Range("J32").Value = 672.5
$ws.Range("L32").Value = 672.5
$ws.Range("N32").Value = -1324.5
$ws.Range("H33").Value = 172.1
$ws.Range("I33").Value = 165.75
$ws.Range("K33").Value = 165.75
$ws.Range("M33").Value = 63.25
$ws.Range("H64").Value = 4442.143
$ws.Range("I64").Value = 2819
$ws.Range("J64").Value = 8500
$ws.Range("K64").Value = 2819
$ws.Range("L64").Value = 8500
$ws.Range("M64").Value = -2571
$ws.Range("N64").Value = -8996
$ws.Range("H67").Value = 4442.143
$ws.Range("I67").Value = 2819
$ws.Range("J67").Value = 8500
$ws.Range("K67").Value = 2819
$ws.Range("L67").Value = 8500
$ws.Range("M67").Value = -1961
$ws.Range("N67").Value = -10216
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H98").Value = 536.64703
$ws.Range("I98").Value = 436.92856
$ws.Range("J98").Value = 1002
$ws.Range("K98").Value = 436.92856
$ws.Range("L98").Value = 1002
$ws.Range("M98").Value = 1061.07144
$ws.Range("N98").Value = -3998
$ws.Range("H99").Value = 158.66667
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490
$ws.Range("H122").Value = 536.64703
$ws.Range("I122").Value = 436.92856
$ws.Range("J122").Value = 1002
$ws.Range("K122").Value = 1310.78568
$ws.Range("L122").Value = 3006
$ws.Range("M122").Value = 1139.21432
$ws.Range("N122").Value = -7906
$ws.Range("H132").Value = 11249.714
$ws.Range("I132").Value = 12202.5
$ws.Range("K132").Value = 36607.5
$ws.Range("M132").Value = -34077.5
$ws.Range("H137").Value = 2517.6365
$ws.Range("I137").Value = 1115.6666
$ws.Range("J137").Value = 4200
$ws.Range("K137").Value = 3346.9998
$ws.Range("L137").Value = 12600
$ws.Range("M137").Value = -796.9998000000001
$ws.Range("N137").Value = -17700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 399.66666
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 399
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 399
$ws.Range("M4").Value = -284
$ws.Range("N4").Value = -631
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 1504.3636
$ws.Range("I122").Value = 943.875
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 2831.625
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -381.625
$ws.Range("N122").Value = -13897
$ws.Range("H132").Value = 1491.0834
$ws.Range("I132").Value = 1239.3
$ws.Range("K132").Value = 3717.9
$ws.Range("M132").Value = -1187.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H62").Value = 58333.168
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 58333.168
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H105").Value = 6163389
$ws.Range("H132").Value = 58330
$ws.Range("J132").Value = 58330
$ws.Range("L132").Value = 58330
$ws.Range("N132").Value = -68450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 262.33334
$ws.Range("I7").Value = 480
$ws.Range("K7").Value = 480
$ws.Range("M7").Value = -367
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864
$ws.Range("H99").Value = 3750
$ws.Range("I99").Value = 5000
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -3502
$ws.Range("N99").Value = -5496
$ws.Range("H126").Value = 3750
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -12440
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H134").Value = 2346.6924
$ws.Range("I134").Value = 2471.4167
$ws.Range("J134").Value = 850
$ws.Range("K134").Value = 7414.250100000001
$ws.Range("L134").Value = 2550
$ws.Range("M134").Value = -4879.250100000001
$ws.Range("N134").Value = -7620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10990.3
$ws.Range("I11").Value = 12282.706
$ws.Range("J11").Value = 3666.6667
$ws.Range("K11").Value = 36848.118
$ws.Range("L11").Value = 11000.0001
$ws.Range("M11").Value = -36708.118
$ws.Range("N11").Value = -11280.0001
$ws.Range("H92").Value = 1588.9166
$ws.Range("J92").Value = 2211.3333
$ws.Range("L92").Value = 6633.999899999999
$ws.Range("N92").Value = -9129.999899999999
$ws.Range("H121").Value = 1937
$ws.Range("I121").Value = 477.66666
$ws.Range("J121").Value = 2666.6667
$ws.Range("K121").Value = 1432.99998
$ws.Range("L121").Value = 8000.000100000001
$ws.Range("M121").Value = -122.9999800000001
$ws.Range("N121").Value = -10620.0001
$ws.Range("H132").Value = 2761.4
$ws.Range("I132").Value = 1741.1428
$ws.Range("K132").Value = 15670.2852
$ws.Range("M132").Value = -13140.2852
$ws.Range("H134").Value = 2142.8572
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 232.13637
$ws.Range("I2").Value = 66.2
$ws.Range("K2").Value = 66.2
$ws.Range("M2").Value = 46.8
$ws.Range("H80").Value = 3499.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3499.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3499.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5495.5
$ws.Range("H83").Value = 3499.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3499.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17497.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -27481.5
$ws.Range("H102").Value = 4444.222
$ws.Range("I102").Value = 3666.3333
$ws.Range("J102").Value = 4833.1665
$ws.Range("K102").Value = 3666.3333
$ws.Range("L102").Value = 4833.1665
$ws.Range("M102").Value = -2044.3333
$ws.Range("N102").Value = -8077.1665
$ws.Range("H122").Value = 1813.3334
$ws.Range("I122").Value = 1059.3334
$ws.Range("K122").Value = 3178.0002
$ws.Range("M122").Value = -728.0001999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1660.7142
$ws.Range("I16").Value = 1681.5454
$ws.Range("J16").Value = 1584.3334
$ws.Range("K16").Value = 1681.5454
$ws.Range("L16").Value = 1584.3334
$ws.Range("M16").Value = -1511.5454
$ws.Range("N16").Value = -1924.3334
$ws.Range("H22").Value = 1022.7692
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 1099.5555
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 1099.5555
$ws.Range("M22").Value = -555
$ws.Range("N22").Value = -1689.5555
$ws.Range("H27").Value = 1022.7692
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 1099.5555
$ws.Range("K27").Value = 850
$ws.Range("L27").Value = 1099.5555
$ws.Range("M27").Value = -743
$ws.Range("N27").Value = -1313.5555
$ws.Range("H61").Value = 166672770
$ws.Range("I61").Value = 250005150
$ws.Range("K61").Value = 250005150
$ws.Range("M61").Value = -250004948
$ws.Range("H113").Value = 166672770
$ws.Range("I113").Value = 250005150
$ws.Range("K113").Value = 250005150
$ws.Range("M113").Value = -250002980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4298.0527
$ws.Range("I126").Value = 2588
$ws.Range("K126").Value = 7764
$ws.Range("M126").Value = -5294
$ws.Range("H132").Value = 2918.3333
$ws.Range("I132").Value = 2902
$ws.Range("K132").Value = 8706
$ws.Range("M132").Value = -6176
$ws.Range("H136").Value = 2412.52
$ws.Range("I136").Value = 1631.1765
$ws.Range("J136").Value = 4072.875
$ws.Range("K136").Value = 4893.529500000001
$ws.Range("L136").Value = 12218.625
$ws.Range("M136").Value = -2343.529500000001
$ws.Range("N136").Value = -17318.625
